$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.552.64'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '1.955.23'
$ws.Range("E3").Value = '  +0.88%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = "'" + '242.91'
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").Value = "'" + '0.625'
$ws.Range("E6").Value = '  +2.78%  '
$ws.Range("D7").Value = "'" + '60.35'
$ws.Range("E7").Value = '  +7.08%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +5.47%  '
$ws.Range("D10").Value = "'" + '0.0790'
$ws.Range("E10").Value = '  -1.55%  '
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("D12").Value = "'" + '14.18'
$ws.Range("E12").Value = '  +6.93%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = "'" + '0.840'
$ws.Range("E13").Value = '  +4.77%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.239.52'
$ws.Range("E14").Value = '  +0.72%  '
$ws.Range("D15").Value = "'" + '21.56'
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("D16").Value = "'" + '5.26'
$ws.Range("E16").Value = '  +2.72%  '
$ws.Range("D17").Value = '1.965.38'
$ws.Range("E17").Value = '  +1.49%  '
$ws.Range("D18").Value = '36.461.84'
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").Value = "'" + '69.25'
$ws.Range("E19").Value = '  +0.65%  '
$ws.Range("D20").Value = '0.0₃0854'
$ws.Range("D21").Value = "'" + '229.37'
$ws.Range("E21").Value = '  +1.26%  '
$ws.Range("D22").Value = "'" + '5.07'
$ws.Range("E22").Value = '  +3.01%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  +2.62%  '
$ws.Range("E25").Value = '  +3.82%  '
$ws.Range("D26").Value = "'" + '0.142'
$ws.Range("E26").Value = '  +8.15%  '
$ws.Range("D27").Value = "'" + '9.16'
$ws.Range("E27").Value = '  +1.08%  '
$ws.Range("D28").Value = "'" + '160.49'
$ws.Range("E28").Value = '  +0.55%  '
$ws.Range("D29").Value = "'" + '19.28'
$ws.Range("E29").Value = '  +1.24%  '
$ws.Range("E30").Value = '  +20.09%  '
$ws.Range("E31").Value = '  +2.42%  '
$ws.Range("E32").Value = '  +4.93%  '
$ws.Range("D33").Value = "'" + '0.0613'
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("D34").Value = "'" + '4.45'
$ws.Range("E34").Value = '  +8.11%  '
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").Value = "'" + '3.47'
$ws.Range("E35").Value = '  +9.44%  '
$ws.Range("B36").Value = 'BinanceUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D36").Value = "'" + '1.00'
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("D37").Value = "'" + '2.27'
$ws.Range("E37").Value = '  +4.14%  '
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("E39").Value = '  -10.53%  '
$ws.Range("E40").Value = '  -2.01%  '
$ws.Range("D41").Value = "'" + '2.91'
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("E42").Value = '  +2.67%  '
$ws.Range("E43").Value = '  +1.33%  '
$ws.Range("D44").Value = "'" + '15.86'
$ws.Range("E44").Value = '  +1.57%  '
$ws.Range("D45").Value = '1.362.19'
$ws.Range("E45").Value = '  +2.53%  '
$ws.Range("D46").Value = "'" + '88.77'
$ws.Range("E46").Value = '  +3.84%  '
$ws.Range("D47").Value = "'" + '1.03'
$ws.Range("E47").Value = '  +0.62%  '
$ws.Range("D48").Value = "'" + '7.16'
$ws.Range("E48").Value = '  +0.73%  '
$ws.Range("E49").Value = '  +0.80%  '
$ws.Range("D50").Value = "'" + '45.11'
$ws.Range("E50").Value = '  +5.17%  '
$ws.Range("D51").Value = '2.135.00'
$ws.Range("E51").Value = '  +1.03%  '
